$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 460.21054
$ws.Range("I33").Value = 446.375
$ws.Range("J33").Value = 534
$ws.Range("K33").Value = 446.375
$ws.Range("L33").Value = 534
$ws.Range("M33").Value = -217.375
$ws.Range("N33").Value = -992

$ws.Range("H40").Value = 2280
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2280
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2280
$ws.Range("N40").Value = -2630
$ws.Range("M40").ClearContents()

$ws.Range("H64").Value = 2334504.8
$ws.Range("I64").Value = 4528244.5
$ws.Range("J64").Value = 3656.25
$ws.Range("K64").Value = 4528244.5
$ws.Range("L64").Value = 3656.25
$ws.Range("M64").Value = -4527996.5
$ws.Range("N64").Value = -4152.25

$ws.Range("H67").Value = 2334504.8
$ws.Range("I67").Value = 4528244.5
$ws.Range("J67").Value = 3656.25
$ws.Range("K67").Value = 4528244.5
$ws.Range("L67").Value = 3656.25
$ws.Range("M67").Value = -4527386.5
$ws.Range("N67").Value = -5372.25

$ws.Range("H74").Value = 3878.879
$ws.Range("I74").Value = 3466.611
$ws.Range("J74").Value = 4373.6
$ws.Range("K74").Value = 3466.611
$ws.Range("L74").Value = 4373.6
$ws.Range("M74").Value = -2530.611
$ws.Range("N74").Value = -6245.6

$ws.Range("H76").Value = 3226.2927
$ws.Range("I76").Value = 3153.125
$ws.Range("K76").Value = 3153.125
$ws.Range("M76").Value = -2838.125

$ws.Range("H77").Value = 3878.879
$ws.Range("I77").Value = 3466.611
$ws.Range("J77").Value = 4373.6
$ws.Range("K77").Value = 17333.055
$ws.Range("L77").Value = 21868
$ws.Range("M77").Value = -12653.055
$ws.Range("N77").Value = -31228

$ws.Range("H79").Value = 3226.2927
$ws.Range("I79").Value = 3153.125
$ws.Range("K79").Value = 3153.125
$ws.Range("M79").Value = -2061.125

$ws.Range("H137").Value = 20050810
$ws.Range("I137").Value = 1048.6666
$ws.Range("J137").Value = 50125450
$ws.Range("K137").Value = 3145.9998
$ws.Range("L137").Value = 150376350
$ws.Range("M137").Value = -595.9998000000001
$ws.Range("N137").Value = -150381450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 781.5
$ws.Range("I2").Value = 744.13635
$ws.Range("J2").Value = 850
$ws.Range("K2").Value = 744.13635
$ws.Range("L2").Value = 850
$ws.Range("M2").Value = -631.13635
$ws.Range("N2").Value = -1076

$ws.Range("H63").Value = 31251506
$ws.Range("I63").Value = 33334714
$ws.Range("K63").Value = 33334714
$ws.Range("M63").Value = -33334028

$ws.Range("H66").Value = 31251506
$ws.Range("I66").Value = 33334714
$ws.Range("K66").Value = 166673570
$ws.Range("M66").Value = -166670138

$ws.Range("H110").Value = 811.9091
$ws.Range("I110").Value = 704.4286
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 704.4286
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1340.5714
$ws.Range("N110").Value = -5090

$ws.Range("H116").Value = 781.5
$ws.Range("I116").Value = 744.13635
$ws.Range("J116").Value = 850
$ws.Range("K116").Value = 744.13635
$ws.Range("L116").Value = 850
$ws.Range("M116").Value = 1549.86365
$ws.Range("N116").Value = -5438

$ws.Range("H132").Value = 7577617
$ws.Range("I132").Value = 8622202
$ws.Range("J132").Value = 4374.5
$ws.Range("K132").Value = 25866606
$ws.Range("L132").Value = 13123.5
$ws.Range("M132").Value = -25864076
$ws.Range("N132").Value = -18183.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 781.5
$ws.Range("I3").Value = 744.13635
$ws.Range("J3").Value = 850
$ws.Range("K3").Value = 744.13635
$ws.Range("L3").Value = 850
$ws.Range("M3").Value = -630.13635
$ws.Range("N3").Value = -1078

$ws.Range("H105").Value = 3749.8462
$ws.Range("I105").Value = 2292.8572
$ws.Range("J105").Value = 4565.76
$ws.Range("K105").Value = 2292.8572
$ws.Range("L105").Value = 4565.76
$ws.Range("M105").Value = -545.8571999999999
$ws.Range("N105").Value = -8059.76

$ws.Range("H134").Value = 4917.8438
$ws.Range("I134").Value = 4394.5386
$ws.Range("J134").Value = 5275.8945
$ws.Range("K134").Value = 13183.6158
$ws.Range("L134").Value = 15827.6835
$ws.Range("M134").Value = -10648.6158
$ws.Range("N134").Value = -20897.6835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3088.25
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3088.25
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3088.25
$ws.Range("N62").Value = -4336.25
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 3088.25
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3088.25
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15441.25
$ws.Range("N65").Value = -21681.25
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8818.299999999999
$ws.Range("I70").Value = 11436.73
$ws.Range("J70").Value = 3955.5
$ws.Range("K70").Value = 11436.73
$ws.Range("L70").Value = 3955.5
$ws.Range("M70").Value = -11166.73
$ws.Range("N70").Value = -4495.5

$ws.Range("H73").Value = 8818.299999999999
$ws.Range("I73").Value = 11436.73
$ws.Range("J73").Value = 3955.5
$ws.Range("K73").Value = 11436.73
$ws.Range("L73").Value = 3955.5
$ws.Range("M73").Value = -10500.73
$ws.Range("N73").Value = -5827.5

$ws.Range("H80").Value = 10253896
$ws.Range("I80").Value = 22224716
$ws.Range("J80").Value = 1703310.4
$ws.Range("K80").Value = 22224716
$ws.Range("L80").Value = 1703310.4
$ws.Range("M80").Value = -22223718
$ws.Range("N80").Value = -1705306.4

$ws.Range("H83").Value = 10253896
$ws.Range("I83").Value = 22224716
$ws.Range("J83").Value = 1703310.4
$ws.Range("K83").Value = 111123580
$ws.Range("L83").Value = 8516552
$ws.Range("M83").Value = -111118588
$ws.Range("N83").Value = -8526536

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7451.579
$ws.Range("J7").Value = 5348.75
$ws.Range("L7").Value = 5348.75
$ws.Range("N7").Value = -5572.75

$ws.Range("H126").Value = 7451.579
$ws.Range("J126").Value = 5348.75
$ws.Range("L126").Value = 16046.25
$ws.Range("N126").Value = -20986.25

$ws.Range("H136").Value = 22734150
$ws.Range("I136").Value = 23811966
$ws.Range("J136").Value = 100005
$ws.Range("K136").Value = 71435898
$ws.Range("L136").Value = 300015
$ws.Range("M136").Value = -71433348
$ws.Range("N136").Value = -305115

$ws.Range("H140").Value = 61900
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 61900
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 61900
$ws.Range("N140").Value = -72260
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1280
$ws.Range("I81").Value = 650
$ws.Range("J81").Value = 1490
$ws.Range("K81").Value = 1300
$ws.Range("L81").Value = 2980
$ws.Range("M81").Value = -239
$ws.Range("N81").Value = -5102

$ws.Range("H84").Value = 1280
$ws.Range("I84").Value = 650
$ws.Range("J84").Value = 1490
$ws.Range("K84").Value = 6500
$ws.Range("L84").Value = 14900
$ws.Range("M84").Value = -1196
$ws.Range("N84").Value = -25508
